# Commit: "Cleaned up code and experimental data"
#
# 1. The "run 1" worksheet (the pre-10-fold-cross-validation experiment) is
#    removed entirely - it was an earlier, less representative experiment.
# 2. The descriptive comment cell on each remaining run sheet is rewritten
#    to number the surviving experiments 1st..5th in order.

$wb = $excel.ActiveWorkbook

# --- 1. Remove the obsolete "run 1" worksheet -------------------------------
$ws = $wb.Worksheets.Item("run 1")
$ws.Delete()

# --- 2. Renumber / reword the experiment-description comments -------------
$ws1 = $wb.Worksheets.Item("bRun 1")
$ws1.Range("D24").Value = "1st successful run of experiment. Result is obtained by randomly shuffling and taking 20% of entire training data to be the validation set."

$ws2 = $wb.Worksheets.Item("run 2")
$ws2.Range("D24").Value = "2nd successful experiment. Uses 10-fold cross validation, which is more representative of this algorithm when classifying SMS"

$ws3 = $wb.Worksheets.Item("run 3")
$ws3.Range("D24").Value = "3rd successful experiment. Uses 10-fold cross validation, which is more representative of this algorithm when classifying SMS"

$ws4 = $wb.Worksheets.Item("run 4")
$ws4.Range("D24").Value = "4th successful experiment. Uses 10-fold cross validation, which is more representative of this algorithm when classifying SMS"

$ws5 = $wb.Worksheets.Item("Run 5")
$ws5.Range("D25").Value = "5th successful experiment. Uses 10-fold cross validation, which is more representative of this algorithm when classifying SMS"
